$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 37; $row++) {
    $cell = $ws.Cells.Item($row, 8)  # Column H
    $v = $cell.Value()
    if ($v -eq $true) {
        $cell.Value = "Accepted"
    } elseif ($v -eq $false) {
        $cell.Value = "Rejected"
    }
}
